# Weekly update: a new week of price data is inserted at the top of the
# data block (row 61), pushing every existing data row down by one. The
# last existing row (129) ends up re-appended as the new last row (130).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 61 - this shifts rows 61:129 down to 62:130,
# carrying their formatting (incl. the date style on column D) with them.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with this week's record. It mirrors
# the old row 61 except for the reporting date (D) and the volume (J).
$ws.Cells.Item(61, 1).Value = 8
$ws.Cells.Item(61, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(61, 3).Value = "Coquimbo"
$ws.Cells.Item(61, 4).Value = 44664
$ws.Cells.Item(61, 5).Value = 4
$ws.Cells.Item(61, 6).Value = 100112001
$ws.Cells.Item(61, 7).Value = "Berenjena"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 520
$ws.Cells.Item(61, 11).Value = 8000
$ws.Cells.Item(61, 12).Value = 9000
$ws.Cells.Item(61, 13).Value = 8500
$ws.Cells.Item(61, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(61, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(61, 16).Value = 170
$ws.Cells.Item(61, 17).Value = 50
$ws.Cells.Item(61, 18).Value = "Hortaliza"
